$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the confidence-interval text for the first data row:
# " (0, 1.29)" -> " (0.00, 1.29)"
$ws.Range("C2").Value = " (0.00, 1.29)"

# Apply a "0.00" number format to the estimate column (B2:B7)
$ws.Range("B2:B7").NumberFormat = "0.00"

# Move the active selection to C7 (last cell of the table)
$ws.Range("C7").Select()
